$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($cell, $text)
    if ($text -match "^[+-]?\d+(\.\d+)?%?$") {
        $cell.Value = "'" + $text
        $cell.ClearFormats()
    } else {
        $cell.Value = $text
    }
}

$data = @(
    @(1, 'HEPARINA SODICA SUBCUT 5000UI', '5000ui', 'CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA', '102980371', 'OK'),
    @(4, 'PROMETAZINA 25MG, CLORIDRATO', '25mg', 'LABORATÓRIO TEUTO BRASILEIRO S/A', '103700321', 'OK'),
    @(5, 'HALOPERIDOL 5MG', '5mg', 'CELLERA FARMACÊUTICA S.A.', '112360011', 'OK'),
    @(6, 'CLORPROMAZINA 40MG/ML SOL ORAL', '40mg/ml', 'SANOFI MEDLEY FARMACÊUTICA LTDA', 'Último registro encontrado: 183260385', 'Pendente'),
    @(7, 'HALOPERIDOL 2MG/ML SOL ORAL', '2mg/ml', 'CELLERA FARMACÊUTICA S.A.', '112360011', 'OK'),
    @(8, 'CLORPROMAZINA 25MG', '25mg', 'SANOFI MEDLEY FARMACÊUTICA LTDA', '183260385', 'OK'),
    @(9, 'CODEINA 30MG', '30mg', 'CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA', '102980199', 'OK'),
    @(10, 'IMIPRAMINA 25MG', '25mg', 'CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA', '102980023', 'OK'),
    @(11, 'RISPERIDONA 3MG', '3mg', 'VITAMEDIC INDUSTRIA FARMACEUTICA LTDA', '103920197', 'OK'),
    @(13, 'RISPERIDONA 1MG', '1mg', 'VITAMEDIC INDUSTRIA FARMACEUTICA LTDA', '103920197', 'OK'),
    @(14, 'LEVOMEPROMAZINA 4% GOTAS', '4%', 'SANOFI MEDLEY FARMACÊUTICA LTDA', 'Último registro encontrado: 183260316', 'Pendente'),
    @(16, 'LIDOCAINA 2% C/ VASO CONSTRITO', '2%', 'CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA', '102980249', 'OK'),
    @(17, 'NITRATO DE CERIO +SULFADIAZINA', 'Concentração não encontrada', 'CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA', 'Último registro encontrado: 102980560', 'Pendente'),
    @(18, 'COLAGENASE+CLORAFENICOL POMADA 30g', '30g', 'CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA', '102980431', 'OK')
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    Set-CellText $ws.Cells.Item($r, 2) $row[1]
    Set-CellText $ws.Cells.Item($r, 3) $row[2]
    Set-CellText $ws.Cells.Item($r, 4) $row[3]
    Set-CellText $ws.Cells.Item($r, 5) $row[4]
    Set-CellText $ws.Cells.Item($r, 6) $row[5]
    $r = $r + 1
}
